$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.057.96"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.23%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.863.46"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  +0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "311.60"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4979"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.23%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3904"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.77%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.09739"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +26.09%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.140"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.55%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "40.87"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "6.463"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.81%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "20.83"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.46%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.862.62"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.45%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.11%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "7.366"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.10%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001124"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +4.77%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "92.91"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06587"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.43"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  +2.49%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.114.35"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.31"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.27%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.289"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.06%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.542"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.51%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.081.37"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.69%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "21.06"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.88%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "156.76"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.38%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "127.47"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.1054"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.89%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.055"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.613"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.17%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.642"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.17%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.06735"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.03%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "9.424"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02392"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.95%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2176"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.38%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.000"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.46%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "11.45"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.27%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6266"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.15%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.175"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("E43").Value = "  +0.00%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "13.48"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.05%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.6013"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.271"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.649"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.77%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "123.96"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.02%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.974"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.86%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.193"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06833"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.50%  "
